# Convert the "commission_type" column (K) from text codes ("p" / "f")
# to a plain integer field (0 / 1), matching commit "change commtype to
# intergerfield".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 11).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 66 }

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 11)
    $txt = $cell.Text
    if ($txt -eq "p") {
        $cell.Value = 0
    } elseif ($txt -eq "f") {
        $cell.Value = 1
    }
}

# Reflect the author's final cursor/scroll position from the edit session.
$ws.Range("J42").Select()
